# Update Betfair Back/Lay odds for 2025-12-25 games.
# Applies the numeric corrections described in the diff to the single
# data worksheet of the workbook (rows 2-9, various odds columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.4
$ws.Range("W2").Value = 1.31
$ws.Range("X2").Value = 13

$ws.Range("F3").Value = 1.9
$ws.Range("G3").Value = 1.99
$ws.Range("I3").Value = 5.1
$ws.Range("K3").Value = 3.9
$ws.Range("P3").Value = 1.8
$ws.Range("Q3").Value = 2.1
$ws.Range("S3").Value = 3.65
$ws.Range("V3").Value = 1.25
$ws.Range("W3").Value = 2

$ws.Range("F4").Value = 2.62
$ws.Range("G4").Value = 3.3
$ws.Range("H4").Value = 2.74
$ws.Range("K4").Value = 3.6
$ws.Range("P4").Value = 1.66
$ws.Range("Q4").Value = 2.02
$ws.Range("V4").Value = 1.39

$ws.Range("F5").Value = 1.49
$ws.Range("G5").Value = 1.57
$ws.Range("H5").Value = 7
$ws.Range("V5").Value = 1.13
$ws.Range("W5").Value = 2.64

$ws.Range("K6").Value = 3.7
$ws.Range("U6").Value = 2.18
$ws.Range("Z6").Value = 16.5
$ws.Range("AB6").Value = 13.5
$ws.Range("AJ6").Value = 60
$ws.Range("AK6").Value = 980

$ws.Range("H7").Value = 8.2
$ws.Range("I7").Value = 9.8
$ws.Range("J7").Value = 4.9
$ws.Range("K7").Value = 5.6
$ws.Range("N7").Value = 4.5
$ws.Range("O7").Value = 1.23
$ws.Range("P7").Value = 2.2
$ws.Range("Q7").Value = 1.67
$ws.Range("S7").Value = 2.72
$ws.Range("T7").Value = 1.94
$ws.Range("W7").Value = 3.05
$ws.Range("X7").Value = 990
$ws.Range("Z7").Value = 100
$ws.Range("AA7").Value = 350
$ws.Range("AB7").Value = 9.2
$ws.Range("AC7").Value = 12
$ws.Range("AE7").Value = 160
$ws.Range("AF7").Value = 9.2
$ws.Range("AG7").Value = 10.5
$ws.Range("AJ7").Value = 12.5
$ws.Range("AK7").Value = 15.5
$ws.Range("AM7").Value = 160
$ws.Range("AN7").Value = 6.2

$ws.Range("F8").Value = 1.62
$ws.Range("G8").Value = 1.76
$ws.Range("H8").Value = 6.4
$ws.Range("I8").Value = 9.6
$ws.Range("J8").Value = 3.4
$ws.Range("N8").Value = 2.56
$ws.Range("O8").Value = 1.51
$ws.Range("Q8").Value = 2.32
$ws.Range("S8").Value = 5
$ws.Range("T8").Value = 2.34
$ws.Range("V8").Value = 1.13
$ws.Range("W8").Value = 2.3

$ws.Range("G9").Value = 2.16
$ws.Range("K9").Value = 4.1
$ws.Range("N9").Value = 3.65
$ws.Range("O9").Value = 1.26
$ws.Range("Q9").Value = 1.78
$ws.Range("T9").Value = 1.67
$ws.Range("U9").Value = 2.06
$ws.Range("W9").Value = 1.86
